$d = $word.ActiveDocument

# Locate the anchor paragraph: "The best number of clusters..."
$anchorText = "The best number of clusters is the one where the WCSS starts reducing gradually."
$startPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text
    $t = $t.TrimEnd([char]13)
    if ($t -eq $anchorText) {
        $startPara = $para
        break
    }
}

if ($startPara -eq $null) {
    throw "Could not locate anchor paragraph 'The best number of clusters...'"
}

# The document always ends with this anchor paragraph followed by two empty
# ListParagraph paragraphs (the second one holding the _GoBack bookmark), then
# the section properties. Replace that whole span (start of anchor .. end of
# document content) with the fully rebuilt OOXML for the new section.
$rangeStart = $startPara.Range.Start
$rangeEnd = $d.Content.End
$full = $d.Range($rangeStart, $rangeEnd)

$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="1080"/>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>The best number of clusters is the one where the WCSS starts reducing gradually.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="8"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">Hierarchical clustering: </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">types include </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>Agglomerative</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> &#8211; bottom up approach and </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>Divisive</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> &#8211; top bottom approach</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="1080"/>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>Dendogram</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> &#8211; </w:t>
  </w:r>
  <w:r>
    <w:t>used to find the optimal number of clusters</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="1080"/>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="1080"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="1080"/>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$null = $full.InsertXML($xmlFrag)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
